# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (column E) data for the two workers (KISSIS ESTHEFANI
# RACERO PELEGRINO / CC 1047413557 and MARCOS ARTURO LOPEZ DE LOS REYES /
# CC 1047440109) is re-sorted: instead of interleaving both workers period
# by period (2309, 2309, 2310, 2310, ...), each worker now owns a
# contiguous block of rows (16-22 and 23-29), each listing periods
# 2403 down to 2309 (most recent first). The "Valor Mora" for period 2403
# also drops from 46400 to 29387 on both blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Worker 1: KISSIS ESTHEFANI RACERO PELEGRINO (CC 1047413557) -> rows 16-22
$worker1Doc    = "1047413557"
$worker1Name   = "KISSIS ESTHEFANI RACERO PELEGRINO"
$worker1Salary = 1160000

# Worker 2: MARCOS ARTURO LOPEZ DE LOS REYES (CC 1047440109) -> rows 23-29
$worker2Doc    = "1047440109"
$worker2Name   = "MARCOS ARTURO LOPEZ DE LOS REYES"
$worker2Salary = 828116

# Periods, most recent first, with their "Valor Mora" (only 2403 changed).
$periods = @(
  @{ Period = "2403"; Mora = 29387 },
  @{ Period = "2402"; Mora = 46400 },
  @{ Period = "2401"; Mora = 46400 },
  @{ Period = "2312"; Mora = 46400 },
  @{ Period = "2311"; Mora = 46400 },
  @{ Period = "2310"; Mora = 46400 },
  @{ Period = "2309"; Mora = 46400 }
)

$startRow = 16

for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("C$row").Value = $worker1Doc
    $ws.Range("D$row").Value = $worker1Name
    $ws.Range("E$row").Value = $periods[$i].Period
    $ws.Range("F$row").Value = $periods[$i].Mora
    $ws.Range("G$row").Value = $worker1Salary
}

$startRow2 = 23
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $startRow2 + $i
    $ws.Range("C$row").Value = $worker2Doc
    $ws.Range("D$row").Value = $worker2Name
    $ws.Range("E$row").Value = $periods[$i].Period
    $ws.Range("F$row").Value = $periods[$i].Mora
    $ws.Range("G$row").Value = $worker2Salary
}
